$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A:F so the long result strings are fully visible (as in the
# authored workbook, where these columns were manually resized / autofit).
$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(2).ColumnWidth = 25.16667
$ws.Columns.Item(3).ColumnWidth = 24.66667
$ws.Columns.Item(4).ColumnWidth = 18.16667
$ws.Columns.Item(5).ColumnWidth = 23.66667
$ws.Columns.Item(6).ColumnWidth = 27.5

$ws.Range("H2").Value = "40.5 + 4.1a"
$ws.Range("H3").Value = "41.7 + 4.4a"
$ws.Range("H4").Value = "45.3 + 4.8b"
$ws.Range("H5").Value = "47.8 + 3.8b"

$ws.Range("J2").Value = "61.3 + 2.8a"
$ws.Range("J3").Value = "64.9 + 2.8bc"
$ws.Range("J4").Value = "65.5 1.9°"
$ws.Range("J5").Value = "63.3 + 1.0b"

[void]$ws.Range("J15").Select()
